$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap coin rows whose rank order changed (B = name, C = link) ---
$tmpB = $ws.Range("B16").Value2
$tmpC = $ws.Range("C16").Value2
$ws.Range("B16").Value2 = $ws.Range("B17").Value2
$ws.Range("C16").Value2 = $ws.Range("C17").Value2
$ws.Range("B17").Value2 = $tmpB
$ws.Range("C17").Value2 = $tmpC

$tmpB = $ws.Range("B24").Value2
$tmpC = $ws.Range("C24").Value2
$ws.Range("B24").Value2 = $ws.Range("B25").Value2
$ws.Range("C24").Value2 = $ws.Range("C25").Value2
$ws.Range("B25").Value2 = $tmpB
$ws.Range("C25").Value2 = $tmpC

$tmpB = $ws.Range("B36").Value2
$tmpC = $ws.Range("C36").Value2
$ws.Range("B36").Value2 = $ws.Range("B37").Value2
$ws.Range("C36").Value2 = $ws.Range("C37").Value2
$ws.Range("B37").Value2 = $tmpB
$ws.Range("C37").Value2 = $tmpC

# --- Update Price (D) and Volume(1h) (E) columns with refreshed data ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.897.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.296.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.90"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.89"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.126"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.863.12"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.66"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.938.85"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.265.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "437.54"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.45"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.514"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.444.95"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.02%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.96"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.40"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.94%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.65"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.20"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.46"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.31"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.98"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.80"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.783.40"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.36"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.09"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0661"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.30"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "321.27"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.56"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0272"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
